$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 105 (sheet ALC)
$ws.Range("H105").Value = 49800
$ws.Range("J105").Value = 49800
$ws.Range("L105").Value = 49800
$ws.Range("N105").Value = -56788

# Row 129 (sheet ALC)
$ws.Range("H129").Value = 1030.6735
$ws.Range("I129").Value = 671.2857
$ws.Range("J129").Value = 1090.5714
$ws.Range("K129").Value = 2013.8571
$ws.Range("L129").Value = 3271.7142
$ws.Range("M129").Value = 2986.1429
$ws.Range("N129").Value = -13271.7142

# Row 132 (sheet ALC)
$ws.Range("H132").Value = 2387.5789
$ws.Range("I132").Value = 2347.4285
$ws.Range("K132").Value = 7042.2855
$ws.Range("M132").Value = -4512.2855

# Row 134 (sheet ALC)
$ws.Range("H134").Value = 41747.5
$ws.Range("J134").Value = 41747.5
$ws.Range("L134").Value = 41747.5
$ws.Range("N134").Value = -51887.5

$ws = $wb.Worksheets.Item("ARM")
# Row 5 (sheet ARM)
$ws.Range("H5").Value = 296
$ws.Range("I5").Value = 295.63635
$ws.Range("K5").Value = 295.63635
$ws.Range("M5").Value = -183.63635

# Row 122 (sheet ARM)
$ws.Range("H122").Value = 2384.5715
$ws.Range("I122").Value = 2438.4
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 7315.200000000001
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -4865.200000000001
$ws.Range("N122").Value = -11650

# Row 133 (sheet ARM)
$ws.Range("H133").Value = 28930.5
$ws.Range("J133").Value = 28930.5
$ws.Range("L133").Value = 28930.5
$ws.Range("N133").Value = -33990.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (sheet BSM)
$ws.Range("H4").Value = 296
$ws.Range("I4").Value = 295.63635
$ws.Range("K4").Value = 295.63635
$ws.Range("M4").Value = -180.63635

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (sheet CRP)
$ws.Range("H58").Value = 3036108.2
$ws.Range("I58").Value = 4786735
$ws.Range("J58").Value = 12298.363
$ws.Range("K58").Value = 4786735
$ws.Range("L58").Value = 12298.363
$ws.Range("M58").Value = -4786532
$ws.Range("N58").Value = -12704.363

# Row 94 (sheet CRP)
$ws.Range("H94").Value = 1206
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1206
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1206
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2108

# Row 99 (sheet CRP)
$ws.Range("H99").Value = 4078
$ws.Range("I99").Value = 3306
$ws.Range("J99").Value = 4850
$ws.Range("K99").Value = 3306
$ws.Range("L99").Value = 4850
$ws.Range("M99").Value = -1808
$ws.Range("N99").Value = -7846

# Row 122 (sheet CRP)
$ws.Range("H122").Value = 8962.695
$ws.Range("I122").Value = 4847.8823
$ws.Range("J122").Value = 20621.334
$ws.Range("K122").Value = 14543.6469
$ws.Range("L122").Value = 61864.00199999999
$ws.Range("M122").Value = -12093.6469
$ws.Range("N122").Value = -66764.00199999999

# Row 126 (sheet CRP)
$ws.Range("H126").Value = 4078
$ws.Range("I126").Value = 3306
$ws.Range("J126").Value = 4850
$ws.Range("K126").Value = 9918
$ws.Range("L126").Value = 14550
$ws.Range("M126").Value = -7448
$ws.Range("N126").Value = -19490

# Row 136 (sheet CRP)
$ws.Range("H136").Value = 3036108.2
$ws.Range("I136").Value = 4786735
$ws.Range("J136").Value = 12298.363
$ws.Range("K136").Value = 14360205
$ws.Range("L136").Value = 36895.089
$ws.Range("M136").Value = -14357655
$ws.Range("N136").Value = -41995.089

$ws = $wb.Worksheets.Item("CUL")
# Row 123 (sheet CUL)
$ws.Range("H123").Value = 2767.353
$ws.Range("J123").Value = 2877.8125
$ws.Range("L123").Value = 8633.4375
$ws.Range("N123").Value = -13533.4375

# Row 126 (sheet CUL)
$ws.Range("H126").Value = 3327.4614
$ws.Range("J126").Value = 3899.7
$ws.Range("L126").Value = 11699.1
$ws.Range("N126").Value = -21579.1

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (sheet GSM)
$ws.Range("H102").Value = 3953.25
$ws.Range("I102").Value = 3757.1428
$ws.Range("J102").Value = 4227.8
$ws.Range("K102").Value = 3757.1428
$ws.Range("L102").Value = 4227.8
$ws.Range("M102").Value = -2135.1428
$ws.Range("N102").Value = -7471.8

# Row 113 (sheet GSM)
$ws.Range("H113").Value = 1790.8334
$ws.Range("I113").Value = 2010.1666
$ws.Range("J113").Value = 1352.1666
$ws.Range("K113").Value = 2010.1666
$ws.Range("L113").Value = 1352.1666
$ws.Range("M113").Value = 159.8334
$ws.Range("N113").Value = -5692.1666

# Row 122 (sheet GSM)
$ws.Range("H122").Value = 7524.6665
$ws.Range("I122").Value = 9532.799999999999
$ws.Range("J122").Value = 4177.778
$ws.Range("K122").Value = 28598.4
$ws.Range("L122").Value = 12533.334
$ws.Range("M122").Value = -26148.4
$ws.Range("N122").Value = -17433.334

# Row 126 (sheet GSM)
$ws.Range("J126").Value = 3619.4546
$ws.Range("L126").Value = 10858.3638
$ws.Range("N126").Value = -15798.3638

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (sheet LTW)
$ws.Range("H7").Value = 4697.25
$ws.Range("I7").Value = 4697.25
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4697.25
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4585.25
$ws.Range("N7").ClearContents()

# Row 9 (sheet LTW)
$ws.Range("H9").Value = 1042.8572
$ws.Range("I9").Value = 383.33334
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 383.33334
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = -159.33334
$ws.Range("N9").Value = -5448

# Row 18 (sheet LTW)
$ws.Range("H18").Value = 9500
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 9000
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = -9828
$ws.Range("N18").Value = -9344

# Row 22 (sheet LTW)
$ws.Range("H22").Value = 616.25
$ws.Range("J22").Value = 350
$ws.Range("L22").Value = 350
$ws.Range("N22").Value = -940

# Row 27 (sheet LTW)
$ws.Range("H27").Value = 616.25
$ws.Range("J27").Value = 350
$ws.Range("L27").Value = 350
$ws.Range("N27").Value = -564

# Row 40 (sheet LTW)
$ws.Range("H40").Value = 3296.652
$ws.Range("I40").Value = 2890.5625
$ws.Range("J40").Value = 4224.857
$ws.Range("K40").Value = 2890.5625
$ws.Range("L40").Value = 4224.857
$ws.Range("M40").Value = -2754.5625
$ws.Range("N40").Value = -4496.857

# Row 68 (sheet LTW)
$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 3500
$ws.Range("K68").Value = 3500
$ws.Range("M68").Value = -2751

# Row 71 (sheet LTW)
$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 3500
$ws.Range("K71").Value = 17500
$ws.Range("M71").Value = -13756

# Row 126 (sheet LTW)
$ws.Range("H126").Value = 4697.25
$ws.Range("I126").Value = 4697.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14091.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11621.75
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (sheet WVR)
$ws.Range("H62").Value = 3651.25
$ws.Range("J62").Value = 3502.5
$ws.Range("L62").Value = 3502.5
$ws.Range("N62").Value = -4750.5

# Row 65 (sheet WVR)
$ws.Range("H65").Value = 3651.25
$ws.Range("J65").Value = 3502.5
$ws.Range("L65").Value = 17512.5
$ws.Range("N65").Value = -23752.5

# Row 122 (sheet WVR)
$ws.Range("H122").Value = 3933.2307
$ws.Range("I122").Value = 1356.9286
$ws.Range("J122").Value = 10491.091
$ws.Range("K122").Value = 4070.7858
$ws.Range("L122").Value = 31473.273
$ws.Range("M122").Value = -1620.7858
$ws.Range("N122").Value = -36373.273

# Row 126 (sheet WVR)
$ws.Range("H126").Value = 1700.7368
$ws.Range("I126").Value = 1753.6
$ws.Range("J126").Value = 1502.5
$ws.Range("K126").Value = 5260.799999999999
$ws.Range("L126").Value = 4507.5
$ws.Range("M126").Value = -2790.799999999999
$ws.Range("N126").Value = -9447.5
